$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5605578392
$ws.Range("C2").Value = -224.73525834
$ws.Range("D2").Value = -225.29581618
$ws.Range("E2").Value = -224.5104056355

$ws.Range("B3").Value = -0.5691767054
$ws.Range("C3").Value = -224.66538499
$ws.Range("D3").Value = -225.2345617
$ws.Range("E3").Value = -224.5104056355

$ws.Range("B4").Value = -0.571565164
$ws.Range("C4").Value = -224.64379284
$ws.Range("D4").Value = -225.215358
$ws.Range("E4").Value = -224.5104056355
